$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# flow10 block (rows 9-13)
$ws.Range("C10").Value = 24.33
$ws.Range("F10").Value = 21.73
$ws.Range("C11").Value = 48.4
$ws.Range("C12").Value = 35.89
$ws.Range("F12").Value = 85.43

# flow15 block (rows 16-20)
$ws.Range("B16").Value = "1.58h"
$ws.Range("C16").Value = 39.04
$ws.Range("C17").Value = 49.09
$ws.Range("D17").Value = 0.638
$ws.Range("F17").Value = 54.91
$ws.Range("C18").Value = 1193.3
$ws.Range("D18").Value = 0.5193
$ws.Range("F18").Value = 40000
$ws.Range("B19").Value = "48.48s"
$ws.Range("C19").Value = 675.07
$ws.Range("D19").Value = 0.2527
$ws.Range("F19").Value = 20000

# Update view: scroll and selection
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$ws.Range("G16").Select()
